# Edit script: insert a new "Bshunt" column into the BUS DATA sheet,
# populate shunt susceptance values, fill in a couple of missing Pmin/Pmax
# values, underline a generator capacity figure, and add an underlined
# blank formatting row below the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BUS DATA")

# --- Insert a new column before the old column O (time-series t1..t24) ---
$ws.Columns("O:O").Insert()

# New column header
$ws.Range("O1").Value = "Bshunt"

# Populate the new Bshunt column (rows 2-31) with 0.01 for every bus
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 15).Value = 0.01
}

# Fill in previously-empty Pmin/Pmax values for the slack/reference buses
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 20
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 20
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 20

# Underline the Rampa value for bus 30
$ws.Range("G31").Font.Underline = 2

# Add an underlined, otherwise empty, formatting row right below the table
$ws.Range("P32:AM32").Font.Underline = 2
$ws.Range("P32:AM32").Value = ""

# Restore view/selection state (best effort - mirrors the author's final
# on-screen position after editing the sheet)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 10
$ws.Range("R32").Select() | Out-Null
